$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used row in column C (the "Förändrad" column)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2() -eq 45181) {
        $cell.Value = 45182
    }
}
